# Add option to use non BAU ZEV minimum by subregion
# Inserts a new variable row ("RZSPbS" / "Required ZEV Sales Percentage by
# Subregion") into the "Key to Variables" sheet, just before the existing
# "SoCDTtiNTY" row (new row 263), pushing all following rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Key to Variables")

# Insert a new blank row at position 263 (existing rows 263+ shift down to 264+)
$ws.Rows.Item(263).Insert()

# Copy the formatting of the "optional" F-column cell (old row 273 / VSbS->VTQaZ
# row, now shifted to row 275 after the insert) onto the new F263 cell so it
# picks up the same "optional" fill/style used elsewhere in the sheet.
$ws.Cells.Item(275, 6).Copy()
$ws.Cells.Item(263, 6).PasteSpecial(-4122)

# Populate the new row's values
$ws.Cells.Item(263, 1).Value = "trans"
$ws.Cells.Item(263, 2).Value = "RZSPbS"
$ws.Cells.Item(263, 3).Value = "Required ZEV Sales Percentage by Subregion"
$ws.Cells.Item(263, 6).Value = "optional"
$ws.Cells.Item(263, 7).Value = "You are modeling non-BAU subregional ZEV sales requirements (as opposed to a national ZEV sales standard)"

# Match the row's wrapped-text height (same as other "optional"/explanatory rows)
$ws.Rows.Item(263).RowHeight = 45

# Update the view: selection moves to the new cell, and the "About" sheet
# becomes the active tab (as in the target workbook).
$ws.Range("G264").Select()

$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Activate()
